$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.897.22'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '3.318.40'
$ws.Range('E3').Value = '  -5.33%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''602.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '''162.03'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.95%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.305.01'
$ws.Range('E8').Value = '  -5.57%  '
$ws.Range('D9').Value = '''0.571'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.38%  '
$ws.Range('D10').Value = '''0.181'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.45%  '
$ws.Range('D11').Value = '''6.54'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -9.60%  '
$ws.Range('D12').Value = '''0.524'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -9.94%  '
$ws.Range('D13').Value = '''41.21'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').Value = '''0.0000253'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -7.51%  '
$ws.Range('D15').Value = '3.856.47'
$ws.Range('E15').Value = '  -5.24%  '
$ws.Range('D16').Value = '68.003.26'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('D17').Value = '''7.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -8.05%  '
$ws.Range('D18').Value = '3.325.52'
$ws.Range('E18').Value = '  -5.17%  '
$ws.Range('D19').Value = '''546.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -10.39%  '
$ws.Range('E20').Value = '  -2.51%  '
$ws.Range('D21').Value = '''15.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -8.31%  '
$ws.Range('D22').Value = '''0.793'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -9.13%  '
$ws.Range('D23').Value = '''8.24'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -8.55%  '
$ws.Range('D24').Value = '''89.13'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -8.49%  '
$ws.Range('D25').Value = '''14.17'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -8.81%  '
$ws.Range('D26').Value = '''3.45'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.83%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -11.21%  '
$ws.Range('D29').Value = '''30.63'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -9.11%  '
$ws.Range('D30').Value = '''8.06'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -9.51%  '
$ws.Range('D31').Value = '''7.22'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -10.05%  '
$ws.Range('E32').Value = '  -7.00%  '
$ws.Range('D33').Value = '''2.65'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -10.73%  '
$ws.Range('D34').Value = '''575.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -10.10%  '
$ws.Range('D35').Value = '''6.10'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -10.58%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').Value = '''55.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.0444'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.73%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.0903'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -8.73%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '''9.78'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -8.55%  '
$ws.Range('D41').Value = '''0.136'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('D42').Value = '''2.82'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -20.82%  '
$ws.Range('D43').Value = '3.032.46'
$ws.Range('E43').Value = '  -9.42%  '
$ws.Range('D44').Value = '''2.63'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -8.62%  '
$ws.Range('D45').Value = '''0.275'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -10.37%  '
$ws.Range('D46').Value = '0.0₃0615'
$ws.Range('E46').Value = '  -16.22%  '
$ws.Range('D47').Value = '''28.53'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -11.22%  '
$ws.Range('D48').Value = '''2.22'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -12.48%  '
$ws.Range('E49').Value = '  -7.28%  '
$ws.Range('D51').Value = '''127.45'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.73%  '
